$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TagFormula($row) {
    $ws.Cells.Item($row, 6).Formula = '="Tags.allTags.push(new Tag(""" & A' + $row + ' & """, " & IF(B' + $row + ' = "", "null", """" & B' + $row + ' & """") & ", """ & C' + $row + ' & """" & IF(D' + $row + ' = "", "",  ", " & D' + $row + ') & "));"'
}

# --- Remove obsolete tag rows (bottom-up so earlier row numbers stay valid) ---
# Row 39: svg
$ws.Rows.Item(39).Delete()
# Row 28: msoffice
$ws.Rows.Item(28).Delete()
# Row 8: ontime
$ws.Rows.Item(8).Delete()

# --- Insert new tag rows (bottom-up so earlier row numbers stay valid) ---

# winserver, before wcf (now at row 44)
$ws.Rows.Item(44).Insert()
$ws.Cells.Item(44, 1).Value = "winserver"
$ws.Cells.Item(44, 3).Value = "Windows Server"
$ws.Cells.Item(44, 4).Value = 9
Set-TagFormula 44

# ssrs, before mvc (now at row 26)
$ws.Rows.Item(26).Insert()
$ws.Cells.Item(26, 1).Value = "ssrs"
$ws.Cells.Item(26, 3).Value = "Microsoft SQL Server Reporting Services"
$ws.Cells.Item(26, 4).Value = 8
Set-TagFormula 26

# highcharts, highmaps, before html (now at row 16)
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(16).Insert()
$ws.Cells.Item(16, 1).Value = "highcharts"
$ws.Cells.Item(16, 3).Value = "Highcharts"
$ws.Cells.Item(16, 4).Value = 6
Set-TagFormula 16
$ws.Cells.Item(17, 1).Value = "highmaps"
$ws.Cells.Item(17, 3).Value = "Highmaps"
$ws.Cells.Item(17, 4).Value = 8
Set-TagFormula 17

# --- Update expertise values that changed ---
$ws.Cells.Item(25, 4).Value = 8   # mantis: 6 -> 8
$ws.Cells.Item(30, 4).Value = 8   # mysql: 7 -> 8
$ws.Cells.Item(40, 4).Value = 6   # tfs: 9 -> 6
$ws.Cells.Item(52, 4).Value = 8   # xsd: 7 -> 8

# --- Update sheet view (scroll position / selection) ---
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("J23").Select()

$wb.Saved = $false
